# Add "Date_of_enrollment" column (new column F) to the raw service uptake
# worksheet, right after the DOB column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F (position 6); everything from F onward shifts right.
$ws.Columns.Item(6).Insert()

# Populate the new header cell.
$ws.Range("F1").Value = "Date_of_enrollment"

# Give the new column its own width (matches the template's other custom
# column widths; closest value reachable through the character-width
# round-trip used by ColumnWidth).
$ws.Columns.Item(6).ColumnWidth = 17

# Restore the active cell/selection to A2, as saved by the author.
$ws.Range("A2").Select()
